$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) / Volume(1h) (E) columns, and the two rows where
# the feed reordered two coins (Binance-PegBSC-USD/Fetch.AI at rows 28+30,
# FirstDigitalUSD/NEARProtocol at rows 37+38), matching the new cryptos feed.
#
# Price values that would otherwise parse as plain numbers (e.g. "0.999")
# are entered with a leading apostrophe so Excel keeps storing them as text,
# exactly like the existing sheet's inline-string cells (e.g. "1.00").

# Row 2
$ws.Range("D2").Value = "61.818.20"
$ws.Range("E2").Value = "  -2.87%  "

# Row 3
$ws.Range("D3").Value = "2.498.09"
$ws.Range("E3").Value = "  -4.07%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'549.94"
$ws.Range("E5").Value = "  -3.91%  "

# Row 6
$ws.Range("D6").Value = "'147.12"
$ws.Range("E6").Value = "  -5.33%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  -4.35%  "

# Row 9
$ws.Range("D9").Value = "2.491.26"
$ws.Range("E9").Value = "  -4.24%  "

# Row 10
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  -9.55%  "

# Row 11
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -1.69%  "

# Row 12
$ws.Range("D12").Value = "'5.36"
$ws.Range("E12").Value = "  -8.09%  "

# Row 13
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  -7.09%  "

# Row 14
$ws.Range("D14").Value = "'26.06"
$ws.Range("E14").Value = "  -7.47%  "

# Row 15
$ws.Range("D15").Value = "2.941.97"
$ws.Range("E15").Value = "  -4.30%  "

# Row 16
$ws.Range("D16").Value = "61.686.87"
$ws.Range("E16").Value = "  -2.78%  "

# Row 17
$ws.Range("D17").Value = "'0.0000163"
$ws.Range("E17").Value = "  -8.50%  "

# Row 18
$ws.Range("D18").Value = "2.480.27"
$ws.Range("E18").Value = "  -4.95%  "

# Row 19
$ws.Range("D19").Value = "'11.18"
$ws.Range("E19").Value = "  -6.74%  "

# Row 20
$ws.Range("D20").Value = "'7.00"
$ws.Range("E20").Value = "  -7.23%  "

# Row 21
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  -8.16%  "

# Row 22
$ws.Range("D22").Value = "'320.46"
$ws.Range("E22").Value = "  -6.51%  "

# Row 23
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").Value = "'63.40"
$ws.Range("E24").Value = "  -5.66%  "

# Row 25
$ws.Range("D25").Value = "'1.75"
$ws.Range("E25").Value = "  -2.20%  "

# Row 26
$ws.Range("D26").Value = "'0.0000103"
$ws.Range("E26").Value = "  -5.14%  "

# Row 27
$ws.Range("D27").Value = "2.606.69"
$ws.Range("E27").Value = "  -5.53%  "

# Row 28
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.49"
$ws.Range("E28").Value = "  -4.42%  "

# Row 29
$ws.Range("D29").Value = "'539.46"
$ws.Range("E29").Value = "  -8.32%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.17%  "

# Row 31
$ws.Range("D31").Value = "'8.36"
$ws.Range("E31").Value = "  -8.57%  "

# Row 32
$ws.Range("D32").Value = "'7.58"
$ws.Range("E32").Value = "  -3.98%  "

# Row 33
$ws.Range("E33").Value = "  -9.03%  "

# Row 34
$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  -8.09%  "

# Row 35
$ws.Range("E35").Value = "  -8.79%  "

# Row 36
$ws.Range("D36").Value = "'5.83"
$ws.Range("E36").Value = "  -10.53%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.84"
$ws.Range("E37").Value = "  -10.79%  "

# Row 38
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39
$ws.Range("D39").Value = "'0.376"
$ws.Range("E39").Value = "  -6.98%  "

# Row 40
$ws.Range("D40").Value = "'18.47"
$ws.Range("E40").Value = "  -6.18%  "

# Row 41
$ws.Range("D41").Value = "'144.30"
$ws.Range("E41").Value = "  -7.13%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("D43").Value = "'1.69"
$ws.Range("E43").Value = "  -9.26%  "

# Row 44
$ws.Range("E44").Value = "  -1.87%  "

# Row 45
$ws.Range("D45").Value = "'2.30"
$ws.Range("E45").Value = "  -8.12%  "

# Row 46
$ws.Range("D46").Value = "'148.62"
$ws.Range("E46").Value = "  -4.91%  "

# Row 47
$ws.Range("D47").Value = "'3.56"
$ws.Range("E47").Value = "  -8.91%  "

# Row 48
$ws.Range("D48").Value = "'20.98"
$ws.Range("E48").Value = "  -9.23%  "

# Row 49
$ws.Range("D49").Value = "'0.0532"
$ws.Range("E49").Value = "  -9.46%  "

# Row 50
$ws.Range("E50").Value = "  -6.59%  "

# Row 51
$ws.Range("D51").Value = "'0.0941"
$ws.Range("E51").Value = "  -6.07%  "
